# updated legacy GSC export data
# The rolling GSC export window advanced by one day: the oldest date row
# (2025-10-10) drops off the top of the "Chart" sheet and every remaining
# row shifts up by one, so the data now ends one row earlier than before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the oldest day's row (row 2, just below the header row) - this
# shifts all subsequent rows up by one, matching a new rolling export.
$ws.Rows("2:2").Delete()
